# Add: Move item in inventory
# Adds a new "Smelt" worksheet (Result / Item smelting recipe table) to the
# workbook, mirroring the formatting already used on the "Craft" sheet, and
# updates the various sheet-view / active-tab bookkeeping that Excel updates
# whenever the user navigates around while making the edit.

$wb = $excel.ActiveWorkbook

$blocks = $wb.Worksheets.Item("Blocks")
$items  = $wb.Worksheets.Item("Items")
$craft  = $wb.Worksheets.Item("Craft")

# --- Visit Items & Craft (matches the selection state left behind in the diff) ---
$items.Activate()
$items.Range("A21").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 21
$win.ScrollColumn = 1
$items.Range("A43").Select()

$craft.Activate()
$craft.Range("A14").Select()

# --- Create the new "Smelt" sheet after "Craft" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$smelt = $wb.Worksheets.Add($null, $lastSheet)
$smelt.Name = "Smelt"
$smelt.Activate()

# Column widths (A:B) - matches the 25.7109375-wide columns used elsewhere
# in the workbook (closest value reachable through ColumnWidth quantization).
$smelt.Range("A1:B1").EntireColumn.ColumnWidth = 24.8776041666666667

# Page margins - matches the margins used on every other sheet in the workbook.
$smelt.PageSetup.LeftMargin = 36.850393728
$smelt.PageSetup.RightMargin = 36.850393728
$smelt.PageSetup.TopMargin = 56.692913399999995
$smelt.PageSetup.BottomMargin = 56.692913399999995
$smelt.PageSetup.HeaderMargin = 22.67716464
$smelt.PageSetup.FooterMargin = 22.67716464

# --- Header row, styled like the other tables' header row (bold/italic white
#     on black, centered) by copying the format from Craft!A1:B1 ---
$smelt.Range("A1").Value = "Result"
$smelt.Range("B1").Value = "Item"
$craft.Range("A1:B1").Copy()
$smelt.Range("A1:B1").PasteSpecial(-4122)

# --- Smelting recipe data rows ---
# Column A = Result item, Column B = source Item.
# Cells whose text refers to a "_block_item" get the same left-aligned style
# (s="4") that's used for block-item names elsewhere in the workbook; we grab
# that exact style from Items!A40 (Grass_block_item) by copying formats.
$items.Range("A40").Copy()

$smelt.Range("A2").Value = "Coal"
$smelt.Range("B2").Value = "Wood_Log_block_item"
$smelt.Range("B2").PasteSpecial(-4122)

$smelt.Range("A3").Value = "Coal"
$smelt.Range("B3").Value = "Coal_Ore_block_item"
$smelt.Range("B3").PasteSpecial(-4122)

$smelt.Range("A4").Value = "Iron"
$smelt.Range("B4").Value = "Iron_Ore_block_item"
$smelt.Range("B4").PasteSpecial(-4122)

$smelt.Range("A5").Value = "Iron"
$smelt.Range("B5").Value = "Raw_Iron"

$smelt.Range("A6").Value = "Gold"
$smelt.Range("B6").Value = "Gold_Ore_block_item"
$smelt.Range("B6").PasteSpecial(-4122)

$smelt.Range("A7").Value = "Gold"
$smelt.Range("B7").Value = "Raw_Gold"

$smelt.Range("A8").Value = "Diamond"
$smelt.Range("B8").Value = "Diamond_Ore_block_item"
$smelt.Range("B8").PasteSpecial(-4122)

$smelt.Range("A9").Value = "Emerald"
$smelt.Range("B9").Value = "Emerald_Ore_block_item"
$smelt.Range("B9").PasteSpecial(-4122)

$smelt.Range("A10").Value = "Cooked_Steak"
$smelt.Range("B10").Value = "Raw_Steak"

$smelt.Range("A11").Value = "Stone_block_item"
$smelt.Range("B11").Value = "Cobblestone_block_item"
$smelt.Range("A11").PasteSpecial(-4122)
$smelt.Range("B11").PasteSpecial(-4122)

# --- Leftover formatted-but-empty rows (45:70), same left-aligned style ---
$items.Range("A40").Copy()
$smelt.Range("A45:B70").PasteSpecial(-4122)

# --- Finally land the selection where the diff shows it ---
$smelt.Range("L19").Select()
